# Simulator full-month coverage: populate Rate/Total for Doug Kinsey's
# 2026-01-14 weekly timesheet export (and the mirrored "Jason Schema" sheet).

$wb = $excel.ActiveWorkbook

$rate = 92

# --- Sheet 1: "Weekly Timesheet" ---
$ws1 = $wb.Worksheets.Item("Weekly Timesheet")

# Rows 2-5 hold the daily entries; column C = Hours, E = Rate, F = Total
for ($r = 2; $r -le 5; $r++) {
    $hours = $ws1.Cells.Item($r, 3).Value2
    $ws1.Cells.Item($r, 5).Value = $rate
    $ws1.Cells.Item($r, 6).Value = $hours * $rate
}

# Recompute the subtotal / grand total rows (Total column F only)
$grandTotal = 0
for ($r = 2; $r -le 5; $r++) {
    $grandTotal = $grandTotal + $ws1.Cells.Item($r, 6).Value2
}

$ws1.Cells.Item(7, 6).Value = $grandTotal    # SUBTOTAL
$ws1.Cells.Item(10, 6).Value = $grandTotal   # HOURLY SUBTOTAL
$ws1.Cells.Item(12, 6).Value = $grandTotal   # GRAND TOTAL

# --- Sheet 2: "Jason Schema" ---
$ws2 = $wb.Worksheets.Item("Jason Schema")

# Rows 2-5 hold the daily entries; column E = Hours, F = Rate, G = Total
for ($r = 2; $r -le 5; $r++) {
    $hours = $ws2.Cells.Item($r, 5).Value2
    $ws2.Cells.Item($r, 6).Value = $rate
    $ws2.Cells.Item($r, 7).Value = $hours * $rate
}
